# Reorder the student rows (A2:B8) and remove the grade columns (C:D) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, in the target order (last name, first name)
$data = @(
    @("Matumona", "Noe"),
    @("Zillig", "Nicolas"),
    @("Sarman", "Dominik"),
    @("Kohler", "Alina"),
    @("Matumona", "Nina"),
    @("asdf", "Marlene"),
    @("Kohler", "Nina")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the grade data in columns C and D for the student rows
# (header row 1 - "Note Exakt" / "Note Gerundet" - is left untouched)
$ws.Range("C2:D8").ClearContents()
